$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.277.20"
$ws.Range("E2").Value = "  +2.79%  "
$ws.Range("D3").Value = "1.785.27"
$ws.Range("E3").Value = "  +3.49%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.008"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  +0.49%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "334.99"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +0.64%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.005"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +0.52%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3784"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  +1.78%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3416"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  +1.52%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "47.97"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  -1.09%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.197"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +1.02%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07429"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +0.16%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.006"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +0.44%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "21.84"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +8.69%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.442"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +0.91%  "
$ws.Range("D15").Value = "1.787.36"
$ws.Range("E15").Value = "  +3.46%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "7.012"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  -0.42%  "
$ws.Range("E17").Value = "  +1.51%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.06650"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  +0.19%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "84.25"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +2.96%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "1.003"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +0.30%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "17.27"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +4.46%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.442"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +4.89%  "
$ws.Range("D23").Value = "27.258.88"
$ws.Range("E23").Value = "  +2.78%  "
$ws.Range("E24").Value = "  -2.10%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.460"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +0.30%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.544"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +6.26%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "21.23"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +9.52%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.447"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +2.56%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "150.27"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -0.40%  "
$ws.Range("D30").Value = "1.991.17"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "132.84"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +1.34%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.054"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -1.38%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.083"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +1.96%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.08633"
$ws.Range("D34").ClearFormats()
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "13.14"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +3.21%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.667"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -1.66%  "
$ws.Range("B37").Value = "TheSandbox"
$ws.Range("C37").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.6849"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +10.41%  "
$ws.Range("B38").Value = "InternetComputer(DFINITY)"
$ws.Range("C38").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "5.401"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +0.86%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.06342"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +2.18%  "
$ws.Range("B40").Value = "FraxShare"
$ws.Range("C40").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "8.774"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +4.19%  "
$ws.Range("B41").Value = "VeChain"
$ws.Range("C41").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.02336"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +0.23%  "
$ws.Range("B42").Value = "Algorand"
$ws.Range("C42").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.2186"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +1.65%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.262"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +3.32%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "14.36"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +0.87%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.005"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +0.50%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.6405"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +6.53%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.841"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -1.70%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.109"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +3.33%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "129.14"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +0.21%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.07174"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +0.12%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "79.04"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +2.64%  "
